# datafile and login TC changes
#
# 1. Add a new "Channel 92" row of test data to the "Channel" sheet.
# 2. Add a new "MyChatList" worksheet (new, currently-empty test-case sheet)
#    after "Channel" and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Append new test data row to the "Channel" sheet ---------------
$channel = $wb.Worksheets.Item("Channel")
$channel.Activate() | Out-Null
$channel.Range("A3").Value = "Channel 92"
$channel.Range("A3").Select() | Out-Null

# --- 2. Add the new "MyChatList" worksheet at the end of the workbook -
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$myChatList = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$myChatList.Name = "MyChatList"

# Give the new sheet the same "couple of blank, taller rows" look the
# author started it off with.
$myChatList.Rows.Item(1).RowHeight = 25.35
$myChatList.Rows.Item(2).RowHeight = 25.35
$myChatList.Columns.Item(1).ColumnWidth = 26.12
$myChatList.Columns.Item(2).ColumnWidth = 34.87

# Make the new sheet the active tab, with A1 selected.
$myChatList.Activate() | Out-Null
$myChatList.Range("A1").Select() | Out-Null
